$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 1342.6666
$ws.Range("I103").Value = 373.2
$ws.Range("J103").Value = 2035.1428
$ws.Range("K103").Value = 1119.6
$ws.Range("L103").Value = 6105.428400000001
$ws.Range("M103").Value = -533.5999999999999
$ws.Range("N103").Value = -7277.428400000001

$ws.Range("H112").Value = 1520.1666
$ws.Range("J112").Value = 1523
$ws.Range("L112").Value = 4569
$ws.Range("N112").Value = -6785

$ws.Range("H113").Value = 3832.8235
$ws.Range("J113").Value = 3195.5
$ws.Range("L113").Value = 3195.5
$ws.Range("N113").Value = -9703.5


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 837.3611
$ws.Range("I2").Value = 822.5333000000001
$ws.Range("J2").Value = 911.5
$ws.Range("K2").Value = 822.5333000000001
$ws.Range("L2").Value = 911.5
$ws.Range("M2").Value = -709.5333000000001
$ws.Range("N2").Value = -1137.5

$ws.Range("H61").Value = 15068.487
$ws.Range("I61").Value = 2551.75
$ws.Range("J61").Value = 62006.25
$ws.Range("K61").Value = 2551.75
$ws.Range("L61").Value = 62006.25
$ws.Range("M61").Value = -2339.75
$ws.Range("N61").Value = -62430.25

$ws.Range("H116").Value = 837.3611
$ws.Range("I116").Value = 822.5333000000001
$ws.Range("J116").Value = 911.5
$ws.Range("K116").Value = 822.5333000000001
$ws.Range("L116").Value = 911.5
$ws.Range("M116").Value = 1471.4667
$ws.Range("N116").Value = -5499.5

$ws.Range("H136").Value = 15068.487
$ws.Range("I136").Value = 2551.75
$ws.Range("J136").Value = 62006.25
$ws.Range("K136").Value = 7655.25
$ws.Range("L136").Value = 186018.75
$ws.Range("M136").Value = -5105.25
$ws.Range("N136").Value = -191118.75


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 837.3611
$ws.Range("I3").Value = 822.5333000000001
$ws.Range("J3").Value = 911.5
$ws.Range("K3").Value = 822.5333000000001
$ws.Range("L3").Value = 911.5
$ws.Range("M3").Value = -708.5333000000001
$ws.Range("N3").Value = -1139.5

$ws.Range("H80").Value = 207.15384
$ws.Range("J80").Value = 206.33333
$ws.Range("L80").Value = 206.33333
$ws.Range("N80").Value = -2202.33333

$ws.Range("H83").Value = 207.15384
$ws.Range("J83").Value = 206.33333
$ws.Range("L83").Value = 1031.66665
$ws.Range("N83").Value = -11015.66665

$ws.Range("H134").Value = 3552.9033
$ws.Range("I134").Value = 2485.72
$ws.Range("J134").Value = 7999.5
$ws.Range("K134").Value = 7457.16
$ws.Range("L134").Value = 23998.5
$ws.Range("M134").Value = -4922.16
$ws.Range("N134").Value = -29068.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 2301.111
$ws.Range("J31").Value = 2828.3809
$ws.Range("K31").Value = 2301.111
$ws.Range("L31").Value = 2828.3809
$ws.Range("M31").Value = -2006.111
$ws.Range("N31").Value = -3418.3809

$ws.Range("I34").Value = 2301.111
$ws.Range("J34").Value = 2828.3809
$ws.Range("K34").Value = 2301.111
$ws.Range("L34").Value = 2828.3809
$ws.Range("M34").Value = -2099.111
$ws.Range("N34").Value = -3232.3809

$ws.Range("H58").Value = 3283.7874
$ws.Range("I58").Value = 2843.4524
$ws.Range("K58").Value = 2843.4524
$ws.Range("M58").Value = -2640.4524

$ws.Range("H86").Value = 4387.55
$ws.Range("I86").Value = 2605.6924
$ws.Range("J86").Value = 7696.7144
$ws.Range("K86").Value = 2605.6924
$ws.Range("L86").Value = 7696.7144
$ws.Range("M86").Value = -1482.6924
$ws.Range("N86").Value = -9942.714400000001

$ws.Range("H89").Value = 4387.55
$ws.Range("I89").Value = 2605.6924
$ws.Range("J89").Value = 7696.7144
$ws.Range("K89").Value = 13028.462
$ws.Range("L89").Value = 38483.572
$ws.Range("M89").Value = -7412.462
$ws.Range("N89").Value = -49715.572

$ws.Range("H107").Value = 26681.527
$ws.Range("I107").Value = 40120.74
$ws.Range("J107").Value = 2904.4614
$ws.Range("K107").Value = 40120.74
$ws.Range("L107").Value = 2904.4614
$ws.Range("M107").Value = -38200.74
$ws.Range("N107").Value = -6744.4614

$ws.Range("H136").Value = 3283.7874
$ws.Range("I136").Value = 2843.4524
$ws.Range("K136").Value = 8530.3572
$ws.Range("M136").Value = -5980.3572


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 791
$ws.Range("I5").Value = 751.25
$ws.Range("J5").Value = 810.875
$ws.Range("K5").Value = 2253.75
$ws.Range("L5").Value = 2432.625
$ws.Range("M5").Value = -2141.75
$ws.Range("N5").Value = -2656.625

$ws.Range("H41").Value = 35
$ws.Range("I41").Value = 50
$ws.Range("J41").Value = 20
$ws.Range("K41").Value = 150
$ws.Range("L41").Value = 60
$ws.Range("M41").Value = 188
$ws.Range("N41").Value = -736

$ws.Range("H113").Value = 1618.2667
$ws.Range("I113").Value = 816.7143
$ws.Range("J113").Value = 2319.625
$ws.Range("K113").Value = 2450.1429
$ws.Range("L113").Value = 6958.875
$ws.Range("M113").Value = -280.1428999999998
$ws.Range("N113").Value = -11298.875

$ws.Range("H121").Value = 11768867
$ws.Range("I121").Value = 10069.25
$ws.Range("J121").Value = 15386959
$ws.Range("K121").Value = 30207.75
$ws.Range("L121").Value = 46160877
$ws.Range("M121").Value = -28897.75
$ws.Range("N121").Value = -46163497

$ws.Range("H135").Value = 791
$ws.Range("I135").Value = 751.25
$ws.Range("J135").Value = 810.875
$ws.Range("K135").Value = 6761.25
$ws.Range("L135").Value = 7297.875
$ws.Range("M135").Value = -4226.25
$ws.Range("N135").Value = -12367.875

$ws.Range("H138").Value = 8874.4
$ws.Range("I138").Value = 7189.75
$ws.Range("K138").Value = 21569.25
$ws.Range("M138").Value = -16429.25


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 113.8
$ws.Range("I2").Value = 206.33333
$ws.Range("J2").Value = 74.14286
$ws.Range("K2").Value = 206.33333
$ws.Range("L2").Value = 74.14286
$ws.Range("M2").Value = -93.33332999999999
$ws.Range("N2").Value = -300.14286

$ws.Range("H113").Value = 4466.269
$ws.Range("I113").Value = 4976
$ws.Range("J113").Value = 3771.182
$ws.Range("K113").Value = 4976
$ws.Range("L113").Value = 3771.182
$ws.Range("M113").Value = -2806
$ws.Range("N113").Value = -8111.182


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2835.8438
$ws.Range("I132").Value = 2466.0981
$ws.Range("J132").Value = 4286.385
$ws.Range("K132").Value = 7398.2943
$ws.Range("L132").Value = 12859.155
$ws.Range("M132").Value = -4868.2943
$ws.Range("N132").Value = -17919.155


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 9999
$ws.Range("I6").Value = 9999
$ws.Range("K6").Value = 9999
$ws.Range("M6").Value = -9884

$ws.Range("H113").Value = 1451.0834
$ws.Range("I113").Value = 1507.3158
$ws.Range("J113").Value = 1237.4
$ws.Range("K113").Value = 4521.9474
$ws.Range("L113").Value = 3712.2
$ws.Range("M113").Value = -2351.9474
$ws.Range("N113").Value = -8052.200000000001

